# Update RPAR_holdings.xlsx model holdings snapshot:
# - bump the "as of" date in the confidential disclaimer note (A18)
# - update Weight (D) and Percent Change (E) values for rows 2-15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unprotect sheet (it carries a password) so values can be written, then
# re-protect with the same password afterward.
$ws.Unprotect("lido")

$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-22 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.05726968181965611
$ws.Range("E2").Value = -0.007813222376329199

$ws.Range("D3").Value = 0.02355192497486343
$ws.Range("E3").Value = -0.004131418453669156

$ws.Range("D4").Value = 0.03138315018993915
$ws.Range("E4").Value = -0.001708428246013582

$ws.Range("D5").Value = 0.03053588721341165
$ws.Range("E5").Value = -0.01317957166392103

$ws.Range("D6").Value = 0.03614106134250143
$ws.Range("E6").Value = -0.01230769230769235

$ws.Range("D7").Value = 0.01882115271451805
$ws.Range("E7").Value = -0.008124076809453484

$ws.Range("D8").Value = 0.004662594067702931
$ws.Range("E8").Value = 0.02427597955706995

$ws.Range("D9").Value = 0.006881363987484337
$ws.Range("E9").Value = 0.003847633705271303

$ws.Range("D10").Value = 0.07085236641211269
$ws.Range("E10").Value = -0.006165919282511312

$ws.Range("D11").Value = 0.07093179731616213
$ws.Range("E11").Value = -0.00615901455767065

$ws.Range("D12").Value = 0.1480380235737685
$ws.Range("E12").Value = 0.004364000572327775

$ws.Range("D13").Value = 0.3867669437702011
$ws.Range("E13").Value = 0.0007855459544383603

$ws.Range("D14").Value = 0.1141640526176785
$ws.Range("E14").Value = 0.004748569341288089

$ws.Range("E15").Value = -0.0008406437345234163

$ws.Protect("lido")
